$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "optional / yellow header" format from an existing header cell (C15)
# onto the new header range N15:AE15, then set the header values + comments.
$ws.Range("C15").Copy() | Out-Null
$ws.Range("N15:AE15").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("N15").Value = 'altitude'
$ws.Range("N15").AddComment('The altitude of the sample is the vertical distance between Earth''s surface above Sea Level and the sampled position in the air.') | Out-Null

$ws.Range("O15").Value = 'biomaterial_provider'
$ws.Range("O15").AddComment('name and address of the lab or PI, or a culture collection identifier') | Out-Null

$ws.Range("P15").Value = 'collected_by'
$ws.Range("P15").AddComment('Name of persons or institute who collected the sample') | Out-Null

$ws.Range("Q15").Value = 'culture_collection'
$ws.Range("Q15").AddComment('Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier') | Out-Null

$ws.Range("R15").Value = 'depth'
$ws.Range("R15").AddComment('Depth is defined as the vertical distance below surface, e.g. for sediment or soil samples depth is measured from sediment or soil surface, respectively. Depth can be reported as an interval for subsurface samples.') | Out-Null

$ws.Range("S15").Value = 'env_biome'
$ws.Range("S15").AddComment('Descriptor of the broad ecological context of a sample. Examples include: desert, taiga or deciduous woodland. FAQ, http://trace.ddbj.nig.ac.jp/biosample/faq_e.html#biome-feature-material EnvO (v 2013-06-14) terms can be found via the link: http://www.environmentontology.org/Browse-EnvO') | Out-Null

$ws.Range("T15").Value = 'genotype'
$ws.Range("T15").AddComment('observed genotype') | Out-Null

$ws.Range("U15").Value = 'host_tissue_sampled'
$ws.Range("U15").AddComment('Type of tissue the initial sample was taken from. Controlled vocabulary, http://bioportal.bioontology.org/ontologies/1005') | Out-Null

$ws.Range("V15").Value = 'identified_by'
$ws.Range("V15").AddComment('name of the taxonomist who identified the specimen') | Out-Null

$ws.Range("W15").Value = 'lab_host'
$ws.Range("W15").AddComment('Scientific name and description of the laboratory host used to propagate the source organism or material from which the sample was obtained, e.g., Escherichia coli DH5a, or Homo sapiens HeLa cells') | Out-Null

$ws.Range("X15").Value = 'lat_lon'
$ws.Range("X15").AddComment('The geographical coordinates of the location where the sample was collected. Specify as decimal degrees latitude and longitude in format "d[d.dddd] N|S d[dd.dddd] W|E", eg, 47.94 N 28.12 W') | Out-Null

$ws.Range("Y15").Value = 'mating_type'

$ws.Range("Z15").Value = 'passage_history'
$ws.Range("Z15").AddComment('Number of passages and passage method') | Out-Null

$ws.Range("AA15").Value = 'samp_size'
$ws.Range("AA15").AddComment('Amount or size of sample (volume, mass or area) that was collected') | Out-Null

$ws.Range("AB15").Value = 'serotype'
$ws.Range("AB15").AddComment('Taxonomy below subspecies; a variety (in bacteria, fungi or virus) usually based on its antigenic properties. Same as serovar and serogroup. e.g. serotype="H1N1" in Influenza A virus CY098518.') | Out-Null

$ws.Range("AC15").Value = 'serovar'
$ws.Range("AC15").AddComment('Taxonomy below subspecies; a variety (in bacteria, fungi or virus) usually based on its antigenic properties. Same as serovar and serotype. Sometimes used as species identifier in bacteria with shaky taxonomy, e.g. Leptospira interrogans serovar Hardjo, http://www.ncbi.nlm.nih.gov/Taxonomy/Browser/wwwtax.cgi?mode=Info&id=176&lvl=3&lin=f&srchmode=3&unlock') | Out-Null

$ws.Range("AD15").Value = 'specimen_voucher'
$ws.Range("AD15").AddComment('Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier') | Out-Null

$ws.Range("AE15").Value = 'temperature'
$ws.Range("AE15").AddComment('temperature of the sample at time of sampling') | Out-Null
